$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 14984.75
$ws.Cells.Item(18, 9).Value = 16542.285
$ws.Cells.Item(18, 11).Value = 16542.285
$ws.Cells.Item(18, 13).Value = -16258.285
$ws.Cells.Item(76, 8).Value = 4699.5
$ws.Cells.Item(76, 10).Value = 4999
$ws.Cells.Item(76, 12).Value = 4999
$ws.Cells.Item(76, 14).Value = -5629
$ws.Cells.Item(79, 8).Value = 4699.5
$ws.Cells.Item(79, 10).Value = 4999
$ws.Cells.Item(79, 12).Value = 4999
$ws.Cells.Item(79, 14).Value = -7183
$ws.Cells.Item(113, 8).Value = 130498.75
$ws.Cells.Item(113, 9).Value = 502495
$ws.Cells.Item(113, 11).Value = 502495
$ws.Cells.Item(113, 13).Value = -499241
$ws.Cells.Item(132, 8).Value = 5291.864
$ws.Cells.Item(132, 9).Value = 2898.3809
$ws.Cells.Item(132, 11).Value = 8695.1427
$ws.Cells.Item(132, 13).Value = -6165.1427
$ws.Cells.Item(137, 8).Value = 1243.8889
$ws.Cells.Item(137, 9).Value = 1243.8889
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 3731.6667
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -1181.6667
$ws.Cells.Item(137, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1169.8334
$ws.Cells.Item(2, 9).Value = 1068.45
$ws.Cells.Item(2, 10).Value = 3197.5
$ws.Cells.Item(2, 11).Value = 1068.45
$ws.Cells.Item(2, 12).Value = 3197.5
$ws.Cells.Item(2, 13).Value = -955.45
$ws.Cells.Item(2, 14).Value = -3423.5
$ws.Cells.Item(37, 8).Value = 24995
$ws.Cells.Item(37, 10).Value = 24995
$ws.Cells.Item(37, 12).Value = 24995
$ws.Cells.Item(37, 14).Value = -25541
$ws.Cells.Item(43, 8).Value = 50000
$ws.Cells.Item(43, 10).Value = 50000
$ws.Cells.Item(43, 12).Value = 50000
$ws.Cells.Item(43, 14).Value = -50626
$ws.Cells.Item(45, 8).Value = 7501.2607
$ws.Cells.Item(45, 9).Value = 11302.637
$ws.Cells.Item(45, 10).Value = 4016.6667
$ws.Cells.Item(45, 11).Value = 11302.637
$ws.Cells.Item(45, 12).Value = 4016.6667
$ws.Cells.Item(45, 13).Value = -10925.637
$ws.Cells.Item(45, 14).Value = -4770.6667
$ws.Cells.Item(61, 8).Value = 2193.7058
$ws.Cells.Item(61, 9).Value = 1986.2
$ws.Cells.Item(61, 11).Value = 1986.2
$ws.Cells.Item(61, 13).Value = -1774.2
$ws.Cells.Item(63, 8).Value = 3799.8
$ws.Cells.Item(63, 9).Value = 2874.75
$ws.Cells.Item(63, 11).Value = 2874.75
$ws.Cells.Item(63, 13).Value = -2188.75
$ws.Cells.Item(66, 8).Value = 3799.8
$ws.Cells.Item(66, 9).Value = 2874.75
$ws.Cells.Item(66, 11).Value = 14373.75
$ws.Cells.Item(66, 13).Value = -10941.75
$ws.Cells.Item(116, 8).Value = 1169.8334
$ws.Cells.Item(116, 9).Value = 1068.45
$ws.Cells.Item(116, 10).Value = 3197.5
$ws.Cells.Item(116, 11).Value = 1068.45
$ws.Cells.Item(116, 12).Value = 3197.5
$ws.Cells.Item(116, 13).Value = 1225.55
$ws.Cells.Item(116, 14).Value = -7785.5
$ws.Cells.Item(132, 8).Value = 5376.25
$ws.Cells.Item(132, 9).Value = 4335
$ws.Cells.Item(132, 11).Value = 13005
$ws.Cells.Item(132, 13).Value = -10475
$ws.Cells.Item(136, 8).Value = 2193.7058
$ws.Cells.Item(136, 9).Value = 1986.2
$ws.Cells.Item(136, 11).Value = 5958.6
$ws.Cells.Item(136, 13).Value = -3408.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1169.8334
$ws.Cells.Item(3, 9).Value = 1068.45
$ws.Cells.Item(3, 10).Value = 3197.5
$ws.Cells.Item(3, 11).Value = 1068.45
$ws.Cells.Item(3, 12).Value = 3197.5
$ws.Cells.Item(3, 13).Value = -954.45
$ws.Cells.Item(3, 14).Value = -3425.5
$ws.Cells.Item(35, 8).Value = 35000
$ws.Cells.Item(35, 10).Value = 35000
$ws.Cells.Item(35, 12).Value = 35000
$ws.Cells.Item(35, 14).Value = -35620
$ws.Cells.Item(94, 8).Value = 824.4666999999999
$ws.Cells.Item(94, 9).Value = 824.0833
$ws.Cells.Item(94, 10).Value = 826
$ws.Cells.Item(94, 11).Value = 824.0833
$ws.Cells.Item(94, 12).Value = 826
$ws.Cells.Item(94, 13).Value = -373.0833
$ws.Cells.Item(94, 14).Value = -1728
$ws.Cells.Item(105, 8).Value = 1999.8
$ws.Cells.Item(105, 9).Value = 1774.75
$ws.Cells.Item(105, 11).Value = 1774.75
$ws.Cells.Item(105, 13).Value = -27.75
$ws.Cells.Item(134, 8).Value = 1893.3529
$ws.Cells.Item(134, 9).Value = 1445.2307
$ws.Cells.Item(134, 11).Value = 4335.6921
$ws.Cells.Item(134, 13).Value = -1800.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5294.6875
$ws.Cells.Item(31, 10).Value = 7587.4287
$ws.Cells.Item(31, 12).Value = 7587.4287
$ws.Cells.Item(31, 14).Value = -8177.4287
$ws.Cells.Item(34, 8).Value = 5294.6875
$ws.Cells.Item(34, 10).Value = 7587.4287
$ws.Cells.Item(34, 12).Value = 7587.4287
$ws.Cells.Item(34, 14).Value = -7991.4287
$ws.Cells.Item(44, 8).Value = 14500
$ws.Cells.Item(44, 10).Value = 19500
$ws.Cells.Item(44, 12).Value = 19500
$ws.Cells.Item(44, 14).Value = -20384
$ws.Cells.Item(122, 8).Value = 74135.86
$ws.Cells.Item(122, 9).Value = 144951
$ws.Cells.Item(122, 10).Value = 3320.7144
$ws.Cells.Item(122, 11).Value = 434853
$ws.Cells.Item(122, 12).Value = 9962.143199999999
$ws.Cells.Item(122, 13).Value = -432403
$ws.Cells.Item(122, 14).Value = -14862.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 151.5
$ws.Cells.Item(23, 10).Value = 146.8
$ws.Cells.Item(23, 12).Value = 440.4
$ws.Cells.Item(23, 14).Value = -910.4000000000001
$ws.Cells.Item(38, 8).Value = 68.42856999999999
$ws.Cells.Item(38, 9).Value = 35.88889
$ws.Cells.Item(38, 11).Value = 107.66667
$ws.Cells.Item(38, 13).Value = 239.33333
$ws.Cells.Item(112, 8).Value = 20000
$ws.Cells.Item(112, 9).Value = 20000
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 11).Value = 60000
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 13).Value = -58892
$ws.Cells.Item(112, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 24430.715
$ws.Cells.Item(43, 9).Value = 5254.25
$ws.Cells.Item(43, 10).Value = 49999.332
$ws.Cells.Item(43, 11).Value = 5254.25
$ws.Cells.Item(43, 12).Value = 49999.332
$ws.Cells.Item(43, 13).Value = -5103.25
$ws.Cells.Item(43, 14).Value = -50301.332
$ws.Cells.Item(46, 8).Value = 42444.848
$ws.Cells.Item(46, 9).Value = 40180
$ws.Cells.Item(46, 10).Value = 49994.332
$ws.Cells.Item(46, 11).Value = 40180
$ws.Cells.Item(46, 12).Value = 49994.332
$ws.Cells.Item(46, 13).Value = -40024
$ws.Cells.Item(46, 14).Value = -50306.332
$ws.Cells.Item(58, 8).Value = 21999
$ws.Cells.Item(58, 10).Value = 21999
$ws.Cells.Item(58, 12).Value = 21999
$ws.Cells.Item(58, 14).Value = -22553
$ws.Cells.Item(132, 8).Value = 4615.619
$ws.Cells.Item(132, 9).Value = 5765.5
$ws.Cells.Item(132, 11).Value = 17296.5
$ws.Cells.Item(132, 13).Value = -14766.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 875.43475
$ws.Cells.Item(16, 10).Value = 774.3333
$ws.Cells.Item(16, 12).Value = 774.3333
$ws.Cells.Item(16, 14).Value = -1114.3333
$ws.Cells.Item(61, 8).Value = 90147.266
$ws.Cells.Item(61, 9).Value = 65289.5
$ws.Cells.Item(61, 11).Value = 65289.5
$ws.Cells.Item(61, 13).Value = -65087.5
$ws.Cells.Item(68, 8).Value = 2000
$ws.Cells.Item(68, 9).Value = 2000
$ws.Cells.Item(68, 10).Value = 2000
$ws.Cells.Item(68, 11).Value = 2000
$ws.Cells.Item(68, 12).Value = 2000
$ws.Cells.Item(68, 13).Value = -1251
$ws.Cells.Item(68, 14).Value = -3498
$ws.Cells.Item(71, 8).Value = 2000
$ws.Cells.Item(71, 9).Value = 2000
$ws.Cells.Item(71, 10).Value = 2000
$ws.Cells.Item(71, 11).Value = 10000
$ws.Cells.Item(71, 12).Value = 10000
$ws.Cells.Item(71, 13).Value = -6256
$ws.Cells.Item(71, 14).Value = -17488
$ws.Cells.Item(93, 8).Value = 22169.21
$ws.Cells.Item(93, 9).Value = 4639.364
$ws.Cells.Item(93, 10).Value = 46272.75
$ws.Cells.Item(93, 11).Value = 4639.364
$ws.Cells.Item(93, 12).Value = 46272.75
$ws.Cells.Item(93, 13).Value = -3391.364
$ws.Cells.Item(93, 14).Value = -48768.75
$ws.Cells.Item(113, 8).Value = 90147.266
$ws.Cells.Item(113, 9).Value = 65289.5
$ws.Cells.Item(113, 11).Value = 65289.5
$ws.Cells.Item(113, 13).Value = -63119.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 25569.154
$ws.Cells.Item(2, 10).Value = 100
$ws.Cells.Item(2, 12).Value = 100
$ws.Cells.Item(2, 14).Value = -324
$ws.Cells.Item(54, 8).Value = 44999.5
$ws.Cells.Item(54, 9).Value = 40000
$ws.Cells.Item(54, 11).Value = 40000
$ws.Cells.Item(54, 13).Value = -39480
$ws.Cells.Item(114, 8).Value = 72203.5
$ws.Cells.Item(114, 10).Value = 72203.5
$ws.Cells.Item(114, 12).Value = 72203.5
$ws.Cells.Item(114, 14).Value = -80881.5
